$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: exception case NNT1028 / New Build Purchase, cleared (FALSE) checkboxes D3:K3
$ws.Range("A3").Value = "NNT1028"
$ws.Range("B3").Value = "New Build Purchase"
$ws.Range("D3:K3").Value = $false

# Row 4: exception case NBT1872 / New Build Purchase, cleared (FALSE) checkboxes D4:K4
$ws.Range("A4").Value = "NBT1872"
$ws.Range("B4").Value = "New Build Purchase"
$ws.Range("D4:K4").Value = $false

# Row 5: exception case NNT1227 / New Build Purchase, cleared (FALSE) checkboxes D5:K5
$ws.Range("A5").Value = "NNT1227"
$ws.Range("B5").Value = "New Build Purchase"
$ws.Range("D5:K5").Value = $false
